{"js": "// The edit removes the redundant word \"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5\" (\"of the constellation\")\n// from the phrase \"\u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a9\u03c1\u03af\u03c9\u03bd\u03b1\", leaving\n// \"\u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd  \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a9\u03c1\u03af\u03c9\u03bd\u03b1\" (note the resulting double space, which\n// matches the author's committed text exactly). This phrase appears 4 times\n// in the document body.\nconst results = context.document.body.search(\"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nresults.items.forEach((result) => {\n  result.insertText(\" \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2\", Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# The edit removes the redundant word \"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5\" (\"of the constellation\")\n# from the phrase \"\u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a9\u03c1\u03af\u03c9\u03bd\u03b1\", leaving\n# \"\u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd  \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a9\u03c1\u03af\u03c9\u03bd\u03b1\" (note the resulting double space, which\n# matches the author's committed text exactly). This phrase appears 4 times\n# in the document body; other, unrelated occurrences of \"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5\"\n# (e.g. \"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u039f\u03c1\u03af\u03c9\u03bd\u03b1\") must stay untouched.\n\n$d = $word.ActiveDocument\n\n$searchText = \"\u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2\"\n$replaceText = \" \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $searchText\n$find.Replacement.Text = $replaceText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
